$wb = $excel.ActiveWorkbook

# --- Sheet "Games": append the newly-completed game (previously the first
#     entry of the "Next" schedule, date 45304 vs ORL) as row 39 ---
$gamesWs = $wb.Worksheets.Item("Games")

$newRow = 39
$lastRow = 38

# Columns A..S in order, skipping K (OppID, handled separately as text)
$cols   = @(1,    2,     3, 4,   5,    6,   7,   8,    9,     10,    12,  13,    14,   15,   16,    17,    18, 19)
$values = @(38, 45304,   4, 112, 97.2, 0.5, 7.2, 23.5, 0.191, 115.2, 100, 0.506, 13.7, 24.4, 0.179, 102.9, 1,  1)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $gamesWs.Cells.Item($newRow, $cols[$i]).Value = $values[$i]
}
# K: OppID (text)
$gamesWs.Cells.Item($newRow, 11).Value = "ORL"

# Match the date formatting used by the rest of column B
$gamesWs.Cells.Item($newRow, 2).NumberFormat = $gamesWs.Cells.Item($lastRow, 2).NumberFormat

# --- Sheet "Next": the game that was just played (row 2, date 45304 vs
#     ORL) drops off the upcoming schedule, so remove it and let the
#     remaining rows shift up ---
$nextWs = $wb.Worksheets.Item("Next")
$nextWs.Rows.Item(2).Delete()
